$wb = $excel.ActiveWorkbook

# --- Sheet 2: "Login & Account Management" ---
$ws = $wb.Worksheets.Item(2)

# Fill in the new "Protect" test case rows (21-25)
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Protect "
$ws.Range("C21").Value = "Unauthorise Access "
$ws.Range("D21").Value = "Click the Logout button if it is login. `nAccess directly : http://ec2-35-161-162-8.us-west-2.compute.amazonaws.com/U-Link/index.html"
$ws.Range("E21").Value = "Page should direct to login page"

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "Protect "
$ws.Range("C22").Value = "Unauthorise Access "
$ws.Range("D22").Value = "Click the Logout button if user is login. `nAccess directly : http://ec2-35-161-162-8.us-west-2.compute.amazonaws.com/U-Link/accountManagement.html"
$ws.Range("E22").Value = "Page should direct to login page"

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "Protect "
$ws.Range("C23").Value = "Unauthorise Access "
$ws.Range("D23").Value = "Click the Logout button if user is login. `nAccess directly : http://ec2-35-161-162-8.us-west-2.compute.amazonaws.com/U-Link/upload.html"
$ws.Range("E23").Value = "Page should direct to login page"

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "Protect "
$ws.Range("C24").Value = "Unauthorise Access "
$ws.Range("D24").Value = "Click the Logout button if user is login. `nAccess directly : http://ec2-35-161-162-8.us-west-2.compute.amazonaws.com/U-link/viewScreenings.html"
$ws.Range("E24").Value = "Page should direct to login page"

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "Protect "
$ws.Range("C25").Value = "Unauthorise Access - Normal User Access Create User Page"
$ws.Range("D25").Value = "Login with Test Case 7`nEmail: ulinkassist_executive@hotmail.com `nPassword:  password!23`nAccess directly : http://ec2-35-161-162-8.us-west-2.compute.amazonaws.com/U-Link/accountManagement.html"
$ws.Range("E25").Value = "Page should direct to user home page"

# Row heights for the new rows
$ws.Rows.Item(21).RowHeight = 48
$ws.Rows.Item(22).RowHeight = 64
$ws.Rows.Item(23).RowHeight = 48
$ws.Rows.Item(24).RowHeight = 64
$ws.Rows.Item(25).RowHeight = 96

# Column B a little wider (closest value this runtime's 1/6-increment
# rounding can reach to the target stored width of 7.6640625)
$ws.Columns.Item(2).ColumnWidth = 6.833333333333333

# --- Sheet 4: "Bootstrap" ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(1).AutoFit()

# Re-activate sheet 2 / select A1:G25 so it ends up the active tab
# (matches target workbook activeTab & tabSelected moving to this sheet)
$ws.Activate() | Out-Null
$ws.Range("A1:G25").Select() | Out-Null
